# Update the division-problem worksheet table with a new set of problems.
# Each data row of the single 5-column table is rewritten in place, cell by
# cell, preserving all existing paragraph/run formatting (rFonts, sz, jc).

$d = $word.ActiveDocument
$t = $d.Tables(1)

# Map of Word-1-indexed table row -> new cell values (left to right).
$rowValues = @{
    1  = @("45÷6=", "24÷7=", "59÷8=", "63÷8=", "87÷6=")
    5  = @("61÷3=", "14÷9=", "32÷2=", "41÷7=", "12÷7=")
    9  = @("77÷2=", "53÷2=", "31÷5=", "70÷2=", "85÷4=")
    13 = @("87÷5=", "28÷5=", "38÷9=", "57÷7=", "16÷8=")
    17 = @("77÷4=", "69÷9=", "66÷3=", "72÷8=", "22÷2=")
}

foreach ($rowIndex in $rowValues.Keys) {
    $values = $rowValues[$rowIndex]
    for ($col = 1; $col -le $values.Length; $col++) {
        $cell = $t.Cell($rowIndex, $col)
        $cell.Range.Text = $values[$col - 1]
    }
}
